$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.764.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.14%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.257.74'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.87%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.85%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.37%  '

$ws.Range('E7').Value = '  +0.18%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.580'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.90%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.250.40'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.92%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.67%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.566'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.84%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.01'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.04%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000266'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.31%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '690.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +10.92%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.802.58'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.31%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.25'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.61%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.158.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.48%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.118'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.98%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.281.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.22%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.34%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.64%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.880'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.52%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.81%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.33'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.50%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.24%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.52%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.78%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '32.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.20%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.33'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.37%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.62'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.76%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '575.73'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.74%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.834.64'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.10%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.69'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.22%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.102'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.14%  '

$ws.Range('E36').Value = '  -0.05%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.09'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.42%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.25'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -10.84%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.127'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.30%  '

$ws.Range('B40').Value = 'ApeXProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.34'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.12%  '

$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.57'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.36%  '

$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '31.46'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.83%  '

$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.05'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.22%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0663'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.21%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.324'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.85%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0403'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.04%  '

$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.01'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.56%  '

$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.126'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.24%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.35'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.78%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.48'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.85%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '128.43'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.68%  '
